$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I7").Value = "sv"
$ws.Range("J7").Value = "Statement-opinion"
$ws.Range("I19").Value = "aa"
$ws.Range("J19").Value = "Agree/Accept"
$ws.Range("I22").Value = "sv"
$ws.Range("J22").Value = "Statement-opinion"
$ws.Range("I29").Value = "sv"
$ws.Range("J29").Value = "Statement-opinion"
$ws.Range("I30").Value = "sv"
$ws.Range("J30").Value = "Statement-opinion"
$ws.Range("I31").Value = "aa"
$ws.Range("J31").Value = "Agree/Accept"
$ws.Range("I43").Value = "sd"
$ws.Range("J43").Value = "Statement-non-opinion"
$ws.Range("I56").Value = "sv"
$ws.Range("J56").Value = "Statement-opinion"
$ws.Range("I57").Value = "sv"
$ws.Range("J57").Value = "Statement-opinion"
$ws.Range("I83").Value = "sd"
$ws.Range("J83").Value = "Statement-non-opinion"
$ws.Range("I84").Value = "sd"
$ws.Range("J84").Value = "Statement-non-opinion"
$ws.Range("I102").Value = "aa"
$ws.Range("J102").Value = "Agree/Accept"
$ws.Range("I104").Value = "sv"
$ws.Range("J104").Value = "Statement-opinion"
$ws.Range("I105").Value = "sd"
$ws.Range("J105").Value = "Statement-non-opinion"
$ws.Range("I116").Value = "b"
$ws.Range("J116").Value = "Acknowledge (Backchannel)"
$ws.Range("I129").Value = "sd"
$ws.Range("J129").Value = "Statement-non-opinion"
$ws.Range("I132").Value = "sd"
$ws.Range("J132").Value = "Statement-non-opinion"
$ws.Range("I137").Value = "aa"
$ws.Range("J137").Value = "Agree/Accept"
$ws.Range("I139").Value = "sd"
$ws.Range("J139").Value = "Statement-non-opinion"
$ws.Range("I141").Value = "aa"
$ws.Range("J141").Value = "Agree/Accept"
$ws.Range("I149").Value = "sd"
$ws.Range("J149").Value = "Statement-non-opinion"
$ws.Range("I151").Value = "sd"
$ws.Range("J151").Value = "Statement-non-opinion"
$ws.Range("I168").Value = "sv"
$ws.Range("J168").Value = "Statement-opinion"
$ws.Range("I179").Value = "sd"
$ws.Range("J179").Value = "Statement-non-opinion"
$ws.Range("I181").Value = "sd"
$ws.Range("J181").Value = "Statement-non-opinion"
$ws.Range("I185").Value = "sd"
$ws.Range("J185").Value = "Statement-non-opinion"
$ws.Range("I193").Value = "sd"
$ws.Range("J193").Value = "Statement-non-opinion"
$ws.Range("I195").Value = "sd"
$ws.Range("J195").Value = "Statement-non-opinion"
$ws.Range("I206").Value = "%"
$ws.Range("J206").Value = "Uninterpretable"
$ws.Range("I210").Value = "sv"
$ws.Range("J210").Value = "Statement-opinion"
$ws.Range("I211").Value = "sd"
$ws.Range("J211").Value = "Statement-non-opinion"
$ws.Range("I214").Value = "sv"
$ws.Range("J214").Value = "Statement-opinion"
$ws.Range("I215").Value = "b"
$ws.Range("J215").Value = "Acknowledge (Backchannel)"
$ws.Range("I219").Value = "sd"
$ws.Range("J219").Value = "Statement-non-opinion"
$ws.Range("I227").Value = "sd"
$ws.Range("J227").Value = "Statement-non-opinion"
$ws.Range("I231").Value = "sd"
$ws.Range("J231").Value = "Statement-non-opinion"
$ws.Range("I233").Value = "sv"
$ws.Range("J233").Value = "Statement-opinion"
$ws.Range("I240").Value = "sd"
$ws.Range("J240").Value = "Statement-non-opinion"
$ws.Range("I244").Value = "ba"
$ws.Range("J244").Value = "Appreciation"
$ws.Range("I247").Value = "sd"
$ws.Range("J247").Value = "Statement-non-opinion"
$ws.Range("I260").Value = "sd"
$ws.Range("J260").Value = "Statement-non-opinion"
$ws.Range("I266").Value = "aa"
$ws.Range("J266").Value = "Agree/Accept"
$ws.Range("I270").Value = "sd"
$ws.Range("J270").Value = "Statement-non-opinion"
$ws.Range("I292").Value = "sv"
$ws.Range("J292").Value = "Statement-opinion"
$ws.Range("I298").Value = "sv"
$ws.Range("J298").Value = "Statement-opinion"
$ws.Range("I299").Value = "sv"
$ws.Range("J299").Value = "Statement-opinion"
$ws.Range("I300").Value = "sv"
$ws.Range("J300").Value = "Statement-opinion"
$ws.Range("I302").Value = "aa"
$ws.Range("J302").Value = "Agree/Accept"
$ws.Range("I303").Value = "sv"
$ws.Range("J303").Value = "Statement-opinion"
$ws.Range("I309").Value = "sd"
$ws.Range("J309").Value = "Statement-non-opinion"
$ws.Range("I314").Value = "%"
$ws.Range("J314").Value = "Uninterpretable"
$ws.Range("I318").Value = "sd"
$ws.Range("J318").Value = "Statement-non-opinion"
$ws.Range("I320").Value = "sv"
$ws.Range("J320").Value = "Statement-opinion"
$ws.Range("I322").Value = "sv"
$ws.Range("J322").Value = "Statement-opinion"
$ws.Range("I324").Value = "sv"
$ws.Range("J324").Value = "Statement-opinion"
$ws.Range("I325").Value = "sd"
$ws.Range("J325").Value = "Statement-non-opinion"
$ws.Range("I330").Value = "sv"
$ws.Range("J330").Value = "Statement-opinion"
$ws.Range("I335").Value = "sv"
$ws.Range("J335").Value = "Statement-opinion"
$ws.Range("I344").Value = "ba"
$ws.Range("J344").Value = "Appreciation"
$ws.Range("I345").Value = "sv"
$ws.Range("J345").Value = "Statement-opinion"
$ws.Range("I350").Value = "sv"
$ws.Range("J350").Value = "Statement-opinion"
$ws.Range("I351").Value = "sd"
$ws.Range("J351").Value = "Statement-non-opinion"
$ws.Range("I356").Value = "sd"
$ws.Range("J356").Value = "Statement-non-opinion"
$ws.Range("I360").Value = "sd"
$ws.Range("J360").Value = "Statement-non-opinion"
$ws.Range("I366").Value = "%"
$ws.Range("J366").Value = "Uninterpretable"
$ws.Range("I386").Value = "sd"
$ws.Range("J386").Value = "Statement-non-opinion"
$ws.Range("I396").Value = "sv"
$ws.Range("J396").Value = "Statement-opinion"
$ws.Range("I408").Value = "sd"
$ws.Range("J408").Value = "Statement-non-opinion"
$ws.Range("I412").Value = "sv"
$ws.Range("J412").Value = "Statement-opinion"
$ws.Range("I448").Value = "sv"
$ws.Range("J448").Value = "Statement-opinion"
$ws.Range("I453").Value = "sd"
$ws.Range("J453").Value = "Statement-non-opinion"
$ws.Range("I457").Value = "aa"
$ws.Range("J457").Value = "Agree/Accept"
$ws.Range("I459").Value = "sd"
$ws.Range("J459").Value = "Statement-non-opinion"
$ws.Range("I477").Value = "sd"
$ws.Range("J477").Value = "Statement-non-opinion"
$ws.Range("I490").Value = "%"
$ws.Range("J490").Value = "Uninterpretable"
$ws.Range("I491").Value = "sd"
$ws.Range("J491").Value = "Statement-non-opinion"
$ws.Range("I493").Value = "aa"
$ws.Range("J493").Value = "Agree/Accept"
$ws.Range("I498").Value = "%"
$ws.Range("J498").Value = "Uninterpretable"
$ws.Range("I502").Value = "ba"
$ws.Range("J502").Value = "Appreciation"
$ws.Range("I505").Value = "sd"
$ws.Range("J505").Value = "Statement-non-opinion"
$ws.Range("I506").Value = "sv"
$ws.Range("J506").Value = "Statement-opinion"
$ws.Range("I517").Value = "qy"
$ws.Range("J517").Value = "Yes-No-Question"
